# The commit swaps the presentation's design theme: the main deck theme
# (ppt/theme/theme1.xml, used by the slide master / all slides) goes from
# the custom "Integral" theme back to the stock default "Office Theme"
# palette (the previous "Integral" colours end up preserved as the
# notes-master's theme instead). Font scheme / format scheme (gradients,
# line styles, effects) are identical between the two themes, so the only
# substantive change is the 12-colour theme colour scheme.
#
# Apply it the way PowerPoint's object model exposes theme colours: via
# Slide.ThemeColorScheme, which reads/writes the <a:clrScheme> of the
# theme shared by every slide (they all hang off the single slide master).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office (default) theme colour scheme, in the standard 12-slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB() isn't available in this host, so colours are passed as the
# packed 0xBBGGRR integer PowerPoint's ColorFormat.RGB uses.
$tcs.Colors(1).RGB  = 0         # dk1      #000000
$tcs.Colors(2).RGB  = 16777215  # lt1      #FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      #44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      #E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  #5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  #ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  #A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  #FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  #4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  #70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    #0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink #954F72
